# ---------------------------------------------------------------
# 1) Rename several StructureName labels (shared strings) that are
#    referenced from the "completed" and "remaining" sheets.
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$wsCompleted = $wb.Worksheets.Item("completed")
$wsRemaining = $wb.Worksheets.Item("remaining")

$labelRenames = @{
    2  = "Construction of Irrigation Inlet"
    3  = "Regulator Re-installation (Rehablitation Sub-project)"
    4  = "Construction of Regulator/ Causeway/Drainage Box Outler (New Haor)"
    5  = "Re-excavation of Khal/River (New Haor)"
    6  = "Re-excavation of Khal/River (Rehab Haor)"
    7  = "Construction of Full Embankment (Rehab Haor)"
    8  = "Construction of Submersible Embankment (Rehab Haor)"
    13 = "Construction of Threshing Floor"
}

foreach ($rowNum in $labelRenames.Keys) {
    $newLabel = $labelRenames[$rowNum]
    $wsCompleted.Cells.Item($rowNum, 2).Value = $newLabel
    $wsRemaining.Cells.Item($rowNum, 2).Value = $newLabel
}

# ---------------------------------------------------------------
# 2) Add two new worksheets: "Current_month" and "Prev_cum",
#    placed after "remaining".
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCurrentMonth = $wb.Worksheets.Add($null, $lastSheet)
$wsCurrentMonth.Name = "Current_month"

$wsPrevCum = $wb.Worksheets.Add($null, $wsCurrentMonth)
$wsPrevCum.Name = "Prev_cum"

# ---------------------------------------------------------------
# 3) Populate headers (row 1, columns B:L) on both new sheets.
# ---------------------------------------------------------------
$headers = @("StructureName", "Unit", "As Per 2nd RDPP", "2014-15", "2015-16", "2016-17", "2017-18", "2018-19", "2019-20", "2020-21", "2021-22")

foreach ($sheet in @($wsCurrentMonth, $wsPrevCum)) {
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $sheet.Cells.Item(1, $i + 2).Value = $headers[$i]
    }
}

# ---------------------------------------------------------------
# 4) Populate data rows 2-13 (columns A:L).
#    Current_month = year-over-year increment of the "completed" sheet.
#    Prev_cum      = previous year's cumulative value (completed shifted right).
# ---------------------------------------------------------------
$currentMonthData = @(
    @(0, "Construction of Irrigation Inlet", "Nos", 116, 0, 0, 0, 0, 17, 28, 50, 21),
    @(1, "Regulator Re-installation (Rehablitation Sub-project)", "Nos", 5, 0, 0, 0, 0, 1, 0, 3, 1),
    @(2, "Construction of Regulator/ Causeway/Drainage Box Outler (New Haor)", "Nos", 112, 0, 0, 2, 17, 41, 8, 30, 14),
    @(3, "Re-excavation of Khal/River (New Haor)", "Km", 337.954, 0, 0, 80, 50, 80, 21, 75, 31.94999999999999),
    @(4, "Re-excavation of Khal/River (Rehab Haor)", "Km", 108.974, 0, 0, 0, 0, 48, 26, 34.97, 0),
    @(5, "Construction of Full Embankment (Rehab Haor)", "Km", 67.11, 0, 0, 0, 0, 44, 16, 7.109999999999999, 0),
    @(6, "Construction of Submersible Embankment (Rehab Haor)", "Km", 61.21, 0, 0, 0, 0, 43, 10, 8.21, 0.10000000000000142),
    @(7, "Construction of Submersible Embankment (New Haor)", "Km", 261.653, 0, 0, 28, 40, 96, 24, 50, 23.649999999999977),
    @(8, "Rehablitation of Regulator (New Haor)", "Nos", 7, 0, 0, 0, 0, 1, 2, 4, 0),
    @(9, "Construction of WMG office", "Nos", 30, 0, 0, 0, 0, 3, 3, 13, 11),
    @(10, "Construction of Threshing Floor", "Nos", 5, 0, 0, 0, 0, 0, 0, 3, 2),
    @(11, "Replecement of Gate (Netrokona & Kishoregonj)", "Nos", 86, 0, 0, 26, 26, 1, 1, 19, 13)
)

$prevCumData = @(
    @(0, "Construction of Irrigation Inlet", "Nos", 116, 0, 0, 0, 0, 0, 17, 45, 95),
    @(1, "Regulator Re-installation (Rehablitation Sub-project)", "Nos", 5, 0, 0, 0, 0, 0, 1, 1, 4),
    @(2, "Construction of Regulator/ Causeway/Drainage Box Outler (New Haor)", "Nos", 112, 0, 0, 0, 2, 19, 60, 68, 98),
    @(3, "Re-excavation of Khal/River (New Haor)", "Km", 337.954, 0, 0, 0, 80, 130, 210, 231, 306),
    @(4, "Re-excavation of Khal/River (Rehab Haor)", "Km", 108.974, 0, 0, 0, 0, 0, 48, 74, 108.97),
    @(5, "Construction of Full Embankment (Rehab Haor)", "Km", 67.11, 0, 0, 0, 0, 0, 44, 60, 67.11),
    @(6, "Construction of Submersible Embankment (Rehab Haor)", "Km", 61.21, 0, 0, 0, 0, 0, 43, 53, 61.21),
    @(7, "Construction of Submersible Embankment (New Haor)", "Km", 261.653, 0, 0, 0, 28, 68, 164, 188, 238),
    @(8, "Rehablitation of Regulator (New Haor)", "Nos", 7, 0, 0, 0, 0, 0, 1, 3, 7),
    @(9, "Construction of WMG office", "Nos", 30, 0, 0, 0, 0, 0, 3, 6, 19),
    @(10, "Construction of Threshing Floor", "Nos", 5, 0, 0, 0, 0, 0, 0, 0, 3),
    @(11, "Replecement of Gate (Netrokona & Kishoregonj)", "Nos", 86, 0, 0, 0, 26, 52, 53, 54, 73)
)

for ($r = 0; $r -lt $currentMonthData.Length; $r++) {
    $row = $currentMonthData[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsCurrentMonth.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

for ($r = 0; $r -lt $prevCumData.Length; $r++) {
    $row = $prevCumData[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsPrevCum.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------
# 5) Match formatting: header row (B1:L1) and id column (A2:A13)
#    use the bold/bordered/centered style from the "completed" sheet.
# ---------------------------------------------------------------
foreach ($sheet in @($wsCurrentMonth, $wsPrevCum)) {
    $wsCompleted.Range("B1:L1").Copy() | Out-Null
    $sheet.Range("B1:L1").PasteSpecial(-4122) | Out-Null
    $wsCompleted.Range("A2:A13").Copy() | Out-Null
    $sheet.Range("A2:A13").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false
